$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column K ("15-jun") with the week's data -----------------
$ws.Cells.Item(1, 11).Value = "15-jun"

$kValues = @{
    2  = 0
    3  = 13.695562417517014
    4  = 18.840491590611677
    5  = 18.759559370813314
    6  = 0
    7  = 6.188497351569211
    8  = 4.5189692517910807
    9  = 10.262801375855235
    10 = 17.172310630797565
    11 = 14.761439822931656
    12 = 0
    13 = 16.337426708654956
    14 = 0
    15 = 0
    16 = 16.481133505248394
    17 = 0
    18 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 11).Value = $kValues[$row]
}

# --- Narrow column K / J down to match the new layout ------------------
# (ColumnWidth is expressed in "characters"; the stored OOXML width is
# quantized to whole pixels, so we pick the COM value whose stored width
# lands on the closest achievable increment to the target 7.5703125.)
$ws.Columns.Item(10).ColumnWidth = 6.666666666666666

# --- Update the selection shown in the sheet view -----------------------
[void]$ws.Range("M7").Select()
